$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BIM INFO")

# The sample/example row (row 4) should no longer ship with pre-filled
# example data. Users only need to fill in the key column + the fields
# they want to update, so clear out all the sample values while leaving
# the cell formatting (styles/borders) intact.
$row4 = $ws.Range("B4:N4")
$row4.Value = $null

# B4 had no explicit formatting before, so it collapses to a completely
# blank cell once its value is cleared; re-touch it so the cell itself
# still exists in the sheet (matching the "key column" placeholder).
$ws.Range("B4").Style = "Normal"

# Reflect the fact that the user left the selection on row 4 (the sample
# row) when the workbook was last saved.
$ws.Range("A4:XFD4").Select()
